$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 91, shifting existing rows 91-166 down to 92-167
$ws.Range("A91").EntireRow.Insert()

# Populate the new row 91 with the new weekly data record
$ws.Range("A91").Value = 10
$ws.Range("B91").Value = "Vega Modelo de Temuco"
$ws.Range("C91").Value = "La Araucanía"
$ws.Range("D91").Value = 45068
$ws.Range("E91").Value = 9
$ws.Range("F91").Value = 100112035
$ws.Range("G91").Value = "Bruselas (repollito)"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 80
$ws.Range("K91").Value = 28000
$ws.Range("L91").Value = 28000
$ws.Range("M91").Value = 28000
$ws.Range("N91").Value = "$/malla 15 kilos"
$ws.Range("O91").Value = "Provincia de Quillota"
$ws.Range("P91").Value = 1867
$ws.Range("Q91").Value = 15
$ws.Range("R91").Value = "Hortaliza"
